# Updated symbol list - refresh crypto "Price" values (column D) on sheet1.
# Values must remain text strings (matching the original inlineStr cells),
# not be auto-converted into numeric cells by Excel's type inference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to Text format before assigning, so Excel keeps the
    # numeric-looking string as text instead of coercing it to a Double.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Restore number format / style so we don't leave a stray "Text" format
    # applied to the cell (matches original, unstyled cells).
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2")  "259.39"
Set-TextValue $ws.Range("D3")  "21.56"
Set-TextValue $ws.Range("D4")  "6.115"
Set-TextValue $ws.Range("D5")  "0.06106"
Set-TextValue $ws.Range("D6")  "3.574"
Set-TextValue $ws.Range("D7")  "6.505"
Set-TextValue $ws.Range("D8")  "1.329"
Set-TextValue $ws.Range("D9")  "0.8226"
Set-TextValue $ws.Range("D10") "0.01321"
Set-TextValue $ws.Range("D11") "0.1601"
Set-TextValue $ws.Range("D12") "0.08124"
Set-TextValue $ws.Range("D13") "0.03531"
Set-TextValue $ws.Range("D14") "0.03194"
Set-TextValue $ws.Range("D15") "0.09203"
Set-TextValue $ws.Range("D16") "3.772"
Set-TextValue $ws.Range("D17") "0.001643"
Set-TextValue $ws.Range("D18") "0.04651"
Set-TextValue $ws.Range("D19") "0.006462"
Set-TextValue $ws.Range("D20") "0.006118"
Set-TextValue $ws.Range("D21") "0.001071"
Set-TextValue $ws.Range("D22") "0.0001504"
Set-TextValue $ws.Range("D23") "3.727"
Set-TextValue $ws.Range("D24") "2.269"
Set-TextValue $ws.Range("D26") "0.1245"
Set-TextValue $ws.Range("D40") "0.04633"
Set-TextValue $ws.Range("D41") "0.006996"
Set-TextValue $ws.Range("D42") "0.003700"
Set-TextValue $ws.Range("D43") "0.1112"
Set-TextValue $ws.Range("D44") "0.01170"
Set-TextValue $ws.Range("D45") "0.00006095"
Set-TextValue $ws.Range("D46") "0.0009934"
Set-TextValue $ws.Range("D47") "0.00000000753"
Set-TextValue $ws.Range("D48") "0.8052"
Set-TextValue $ws.Range("D49") "0.001130"
